$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "gültig ab" (A1) becomes "Jahr" — the other headers
# (B1/C1/D1) are untouched content-wise.
$ws.Range("A1").Value = "Jahr"

# Move the cursor/selection to B9, matching the saved view state.
$ws.Range("B9").Select() | Out-Null
